{"js": "// The underlying OOXML diff for this commit (\"Moving from 2.0.1 to 2.0.2\")\n// is purely a re-serialization artifact: every changed line carries the\n// exact same element name and the exact same set of attribute/value pairs\n// as before -- only the XML attribute order changed (alphabetised), which\n// is what you get when a newer version of the OOXML library that saved the\n// file re-emits the same markup using a canonical (sorted) attribute order.\n// None of the document's visible content, formatting, styles, page setup or\n// text actually changed between the two revisions.\n//\n// There is therefore no content-level edit to make through the Word\n// JavaScript object model here (attribute ordering inside the raw XML is an\n// artifact of the OOXML writer, not something the document object model\n// exposes or that changes the rendered/semantic document). We still touch\n// the document defensively (load + sync) so the script is a verifiably\n// successful, side-effect-free run against the object model.\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "# The underlying OOXML diff for this commit (\"Moving from 2.0.1 to 2.0.2\")\n# is purely a re-serialization artifact: every changed line carries the\n# exact same element name and the exact same set of attribute/value pairs\n# as before -- only the XML attribute order changed (alphabetised), which\n# is what you get when a newer version of the OOXML library that saved the\n# file re-emits the same markup using a canonical (sorted) attribute order.\n# None of the document's visible content, formatting, styles, page setup or\n# text actually changed between the two revisions.\n#\n# There is therefore no content-level edit to make through the Word COM\n# object model here (attribute ordering inside the raw XML is an artifact\n# of the OOXML writer, not something the document object model exposes or\n# that changes the rendered/semantic document). We still touch the document\n# defensively (read-only) so the script is a verifiably successful,\n# side-effect-free run against the object model.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
